# Apply the edits described by the commit:
#  - Update three circuit ratings (E2, E3) on "Bus connections" from 7500 -> 10000
#  - Update circuit rating (E9) on "Bus connections" from 15000 -> 20000
#  - Update generator capacities (E4, E5) on "Generator data" from 15000 -> 20000
#  - Downstream formulas recalc automatically
#  - Update the active selection on "Bus connections" and "Generator data"
#  - Make "Bus connections" the active (selected) sheet/tab

$wb = $excel.ActiveWorkbook

$busConnections = $wb.Worksheets.Item("Bus connections")
$generatorData  = $wb.Worksheets.Item("Generator data")

# --- Value edits ---
$busConnections.Range("E2").Value = 10000
$busConnections.Range("E3").Value = 10000
$busConnections.Range("E9").Value = 20000

$generatorData.Range("E4").Value = 20000
$generatorData.Range("E5").Value = 20000

# --- Selection / active sheet updates ---
# Update the selection remembered on "Generator data" first (it ends up NOT
# being the active tab in the final workbook state).
$generatorData.Activate() | Out-Null
$generatorData.Range("E6").Select() | Out-Null

# Finally activate "Bus connections" with its new selection; this becomes
# the workbook's active tab when saved.
$busConnections.Activate() | Out-Null
$busConnections.Range("E4").Select() | Out-Null
